# One-to-many Relationships: Updated Figures to change A$B notation to A.B.
#
# This script:
#  1. Refreshes the cached "datetimeFigureOut" field text (7/6/12 -> 6/15/13)
#     on every slide layout, the slide master and the notes master.
#  2. Renames the "LinkedList$Entry" label on slide 1 to "LinkedList.Entry"
#     and shrinks its textbox to the new (narrower) auto-fit width.

$p = $ppt.ActivePresentation

$oldDate = "7/6/12"
$newDate = "6/15/13"

# --- 1a. Slide master's own "Date Placeholder 3" -------------------------
$master = $p.SlideMaster
$masterDateShape = $master.Shapes.Item("Date Placeholder 3")
if ($masterDateShape.TextFrame.TextRange.Text -eq $oldDate) {
    $masterDateShape.TextFrame.TextRange.Text = $newDate
}

# --- 1b. Every slide layout's "Date Placeholder 3" ------------------------
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    $dateShape = $layout.Shapes.Item("Date Placeholder 3")
    if ($dateShape.TextFrame.TextRange.Text -eq $oldDate) {
        $dateShape.TextFrame.TextRange.Text = $newDate
    }
}

# --- 1c. Notes master's "Date Placeholder 2" -------------------------------
$notesMaster = $p.NotesMaster
$notesDateShape = $notesMaster.Shapes.Item("Date Placeholder 2")
if ($notesDateShape.TextFrame.TextRange.Text -eq $oldDate) {
    $notesDateShape.TextFrame.TextRange.Text = $newDate
}

# --- 2. "LinkedList$Entry" -> "LinkedList.Entry" on slide 1 ---------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "LinkedList`$Entry") {
            $shp.TextFrame.TextRange.Text = "LinkedList.Entry"
            # Shrink the auto-fit textbox to match the shorter rendered text
            # (1727907 EMU -> 1672253 EMU @ 12700 EMU/pt); height (369332 EMU)
            # is unchanged.
            $shp.Width = 131.6735
        }
    }
}
